$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (target stored width 16.28988764044944 chars;
# ColumnWidth is stored internally on a whole-pixel grid at 7 px/char + 5 px
# padding, so request the pixel-exact value that resolves to that width)
$ws.Columns.Item(1).ColumnWidth = 15.575601926163726

# Update row 2 values
$ws.Range("A2").Value = 32145698741
$ws.Range("B2").Value = "Priyanka Muddana"
$ws.Range("C2").Value = 9278.96
$ws.Range("D2").Value = "November"

# Add new row 3
$ws.Range("A3").Value = 123654789963
$ws.Range("B3").Value = "Vidya Sagar pogiri"
$ws.Range("C3").Value = 9793.33
$ws.Range("D3").Value = "November"
